# issue #5: add legislator_id, name, date into dataframe
#
# Adds three new trailing columns - date, legislator_name, legislator_id -
# to the "股票" (stocks) worksheet, filled in for the header row and every
# existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row

# --- Header row: copy the formatting of the existing last header cell
#     (G1, the bold/bordered header style) into H1:J1, then label them. ---
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("G1").Copy($ws.Range("I1"))
$ws.Range("G1").Copy($ws.Range("J1"))

$ws.Cells.Item(1, 8).Value  = "date"
$ws.Cells.Item(1, 9).Value  = "legislator_name"
$ws.Cells.Item(1, 10).Value = "legislator_id"

# --- Data rows ---
for ($r = 2; $r -le $lastRow; $r++) {

    # date: use a leading apostrophe so the "2011-11-21" text isn't
    # auto-converted into a date serial by the smart-entry parser, then
    # strip the resulting quote-prefix formatting back off so the cell
    # ends up as plain text (matching the source data's shared-string cell).
    $ws.Cells.Item($r, 8).Value = "'2011-11-21"
    $ws.Cells.Item($r, 7).Copy()
    $ws.Cells.Item($r, 8).PasteSpecial(-4122)

    # legislator_name / legislator_id: plain text / number, no special
    # handling required.
    $ws.Cells.Item($r, 7).Copy($ws.Cells.Item($r, 9))
    $ws.Cells.Item($r, 9).Value = "羅淑蕾"

    $ws.Cells.Item($r, 7).Copy($ws.Cells.Item($r, 10))
    $ws.Cells.Item($r, 10).Value = 1638
}

$excel.CutCopyMode = $false
